$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr0 = New-Object "object[,]" 24,2
$arr0[0,0] = 2.421508577308259
$arr0[0,1] = 0.3410227133981607
$arr0[1,0] = 2.283797213877733
$arr0[1,1] = 0.3026433378636568
$arr0[2,0] = 2.200732835608108
$arr0[2,1] = 0.279224936085825
$arr0[3,0] = 2.167255407999221
$arr0[3,1] = 0.2697173617881958
$arr0[4,0] = 2.161718888988219
$arr0[4,1] = 0.268140750490943
$arr0[5,0] = 2.200279844996942
$arr0[5,1] = 0.2790965710823343
$arr0[6,0] = 2.373714421932107
$arr0[6,1] = 0.3277582781830688
$arr0[7,0] = 2.725803725058029
$arr0[7,1] = 0.4244101770302677
$arr0[8,0] = 2.992036198106291
$arr0[8,1] = 0.4962620210480964
$arr0[9,0] = 3.114848387848895
$arr0[9,1] = 0.5291530557911983
$arr0[10,0] = 3.161602930466415
$arr0[10,1] = 0.5416391420761784
$arr0[11,0] = 3.151522405849505
$arr0[11,1] = 0.5389486394908545
$arr0[12,0] = 3.118689916745723
$arr0[12,1] = 0.5301796646521666
$arr0[13,0] = 3.098611502498045
$arr0[13,1] = 0.5248124903463349
$arr0[14,0] = 2.984044674942197
$arr0[14,1] = 0.494116768328297
$arr0[15,0] = 2.914199967029219
$arr0[15,1] = 0.4753395340571842
$arr0[16,0] = 2.87418709427476
$arr0[16,1] = 0.46455859851892
$arr0[17,0] = 2.860666789555978
$arr0[17,1] = 0.4609116185696962
$arr0[18,0] = 2.921618475310538
$arr0[18,1] = 0.4773364017517565
$arr0[19,0] = 3.128326849163955
$arr0[19,1] = 0.5327544730786258
$arr0[20,0] = 3.264871306692612
$arr0[20,1] = 0.5691545651745287
$arr0[21,0] = 3.191861225063633
$arr0[21,1] = 0.5497100883060853
$arr0[22,0] = 2.918264126368229
$arr0[22,1] = 0.4764335732895688
$arr0[23,0] = 2.629246040536316
$arr0[23,1] = 0.3981223416926696
$ws.Range("B2:C25").Value = $arr0

$arr1 = New-Object "object[,]" 24,3
$arr1[0,0] = 0.03259823985588861
$arr1[0,1] = 0.4443680307746263
$arr1[0,2] = 0.002523432179149924
$arr1[1,0] = 0.03264722439960344
$arr1[1,1] = 0.387822817061874
$arr1[1,2] = 0.002529511640271156
$arr1[2,0] = 0.03268112333721596
$arr1[2,1] = 0.3531389305168915
$arr1[2,2] = 0.002533437071178822
$arr1[3,0] = 0.03269589889003588
$arr1[3,1] = 0.3390132514313251
$arr1[3,2] = 0.002535085328995271
$arr1[4,0] = 0.03269841044109612
$arr1[4,1] = 0.336668177824194
$arr1[4,2] = 0.002535361962317058
$arr1[5,0] = 0.03268131871206448
$arr1[5,1] = 0.3529483938368969
$arr1[5,2] = 0.002533459102974533
$arr1[6,0] = 0.03261433658489588
$arr1[6,1] = 0.4248636149813336
$arr1[6,2] = 0.002525488507106519
$arr1[7,0] = 0.03251330666946672
$arr1[7,1] = 0.5661985755042025
$arr1[7,2] = 0.002511378390623326
$arr1[8,0] = 0.03245757029810314
$arr1[8,1] = 0.6702781546542269
$arr1[8,2] = 0.002501927021893704
$arr1[9,0] = 0.03243623338993462
$arr1[9,1] = 0.7176906081379002
$arr1[9,2] = 0.002497823659394285
$arr1[10,0] = 0.03242873187049811
$arr1[10,1] = 0.7356546913071611
$arr1[10,2] = 0.002496297835712113
$arr1[11,0] = 0.03243032172674543
$arr1[11,1] = 0.7317853510981394
$arr1[11,2] = 0.002496625205216506
$arr1[12,0] = 0.0324356046417964
$arr1[12,1] = 0.7191683204515869
$arr1[12,2] = 0.002497697568153323
$arr1[13,0] = 0.03243891591028714
$arr1[13,1] = 0.7114413442032514
$arr1[13,2] = 0.002498358066843592
$arr1[14,0] = 0.03245904560812662
$arr1[14,1] = 0.6671810134426437
$arr1[14,2] = 0.002502199118077666
$arr1[15,0] = 0.03247242393957972
$arr1[15,1] = 0.6400460337215605
$arr1[15,2] = 0.002504605586357159
$arr1[16,0] = 0.03248049694043398
$arr1[16,1] = 0.6244449056556647
$arr1[16,2] = 0.002506008192604066
$arr1[17,0] = 0.0324832952512546
$arr1[17,1] = 0.6191636801734006
$arr1[17,2] = 0.002506486268134545
$arr1[18,0] = 0.0324709606540321
$arr1[18,1] = 0.6429339538360921
$arr1[18,2] = 0.002504347503546138
$arr1[19,0] = 0.03243403722156191
$arr1[19,1] = 0.7228739723492197
$arr1[19,2] = 0.002497381829613459
$arr1[20,0] = 0.03241327676887762
$arr1[20,1] = 0.7751780083420101
$arr1[20,2] = 0.002492992664010947
$arr1[21,0] = 0.03242404834099322
$arr1[21,1] = 0.7472568307916134
$arr1[21,2] = 0.002495320358134937
$arr1[22,0] = 0.03247162101673751
$arr1[22,1] = 0.6416283278902171
$arr1[22,2] = 0.00250446412326645
$arr1[23,0] = 0.03253739165902214
$arr1[23,1] = 0.5279251897347308
$arr1[23,2] = 0.002515033985441636
$ws.Range("E2:G25").Value = $arr1

$arr2 = New-Object "object[,]" 24,1
$arr2[0,0] = 1.530761831882074
$arr2[1,0] = 1.522049614036291
$arr2[2,0] = 1.517673056291159
$arr2[3,0] = 1.51613274102418
$arr2[4,0] = 1.515891617605348
$arr2[5,0] = 1.517651300375896
$arr2[6,0] = 1.527555077468094
$arr2[7,0] = 1.554763596097814
$arr2[8,0] = 1.579601132760402
$arr2[9,0] = 1.591975194989061
$arr2[10,0] = 1.596817331255806
$arr2[11,0] = 1.595767509373275
$arr2[12,0] = 1.592370417099914
$arr2[13,0] = 1.590310012402512
$arr2[14,0] = 1.578814243261505
$arr2[15,0] = 1.572038593394851
$arr2[16,0] = 1.568242530939273
$arr2[17,0] = 1.566974563713728
$arr2[18,0] = 1.572749395839239
$arr2[19,0] = 1.5933639684187
$arr2[20,0] = 1.607748911489978
$arr2[21,0] = 1.5999873630698
$arr2[22,0] = 1.572427732733786
$arr2[23,0] = 1.546559035530635
$ws.Range("I2:I25").Value = $arr2

$arr3 = New-Object "object[,]" 24,1
$arr3[0,0] = 0.2966162175013238
$arr3[1,0] = 0.2861052814412375
$arr3[2,0] = 0.2798419397186649
$arr3[3,0] = 0.2773371735172816
$arr3[4,0] = 0.276924125485607
$arr3[5,0] = 0.2798079672212879
$arr3[6,0] = 0.292952348904322
$arr3[7,0] = 0.3202540559151714
$arr3[8,0] = 0.341265843412458
$arr3[9,0] = 0.351036897366356
$arr3[10,0] = 0.3547679119350136
$arr3[11,0] = 0.3539629900692205
$arr3[12,0] = 0.3513432283253621
$arr3[13,0] = 0.3497425867763013
$arr3[14,0] = 0.3406315911155673
$arr3[15,0] = 0.3350970070643342
$arr3[16,0] = 0.3319336586790911
$arr3[17,0] = 0.3308660275759507
$arr3[18,0] = 0.3356841008000657
$arr3[19,0] = 0.352111873678794
$arr3[20,0] = 0.3630288134275048
$arr3[21,0] = 0.3571856110162344
$arr3[22,0] = 0.335418617944768
$arr3[23,0] = 0.3127025252126572
$ws.Range("L2:L25").Value = $arr3

$arr4 = New-Object "object[,]" 24,1
$arr4[0,0] = 1.728444377538857
$arr4[1,0] = 1.748483226842474
$arr4[2,0] = 1.761468412562685
$arr4[3,0] = 1.766930694600397
$arr4[4,0] = 1.767847995829911
$arr4[5,0] = 1.761541388368016
$arr4[6,0] = 1.735212022799978
$arr4[7,0] = 1.689011640221253
$arr4[8,0] = 1.658414317066153
$arr4[9,0] = 1.645229047336919
$arr4[10,0] = 1.640342310560897
$arr4[11,0] = 1.641390019720951
$arr4[12,0] = 1.644824879055413
$arr4[13,0] = 1.646942687384708
$arr4[14,0] = 1.659290831041375
$arr4[15,0] = 1.667054438708597
$arr4[16,0] = 1.671588865529351
$arr4[17,0] = 1.673135973651092
$arr4[18,0] = 1.666220842293626
$arr4[19,0] = 1.643813088297911
$arr4[20,0] = 1.629787922278268
$arr4[21,0] = 1.637216473928831
$arr4[22,0] = 1.66659749006147
$arr4[23,0] = 1.700925028598874
$ws.Range("N2:N25").Value = $arr4

